# "Thumbnails mit absoluten Pfaden"
# Prefix every thumbnail filename (column D, rows 2-15) with the absolute
# base URL so the sheet stores full links instead of bare file names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseUrl = "https://biologie-lernprogramme.de/vorschaubilder/"

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Range("D$row")
    $current = $cell.Value()
    if ($current -and -not ($current.ToString().StartsWith($baseUrl))) {
        $cell.Value = $baseUrl + $current
    }
}

# Column width adjustments: column B got narrower, and the old merged
# C:D width got split into a dedicated column C and a wider column D.
$ws.Columns.Item(2).ColumnWidth = 44.5
$ws.Columns.Item(3).ColumnWidth = 31.333333333333332
$ws.Columns.Item(4).ColumnWidth = 101.5

# Restore the active selection to the cell that was selected when the
# workbook was last saved.
$ws.Activate() | Out-Null
$ws.Range("D21").Select() | Out-Null
